$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 269-388: Fecha (D), Calidad (I), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M), Precio $/Kg (P) shift down
# by two rows (one date block = Primera+Segunda), with a new date block
# inserted at the top (rows 269-270) and the last two source rows pushed
# into two brand-new rows (389-390) at the bottom.
$data = @(
    @(269,44466,"Primera",3200,600,700,650,650),
    @(270,44466,"Segunda",1600,500,550,525,525),
    @(271,44389,"Primera",3160,600,700,650,650),
    @(272,44389,"Segunda",1600,500,550,525,525),
    @(273,44249,"Primera",3000,850,900,875,875),
    @(274,44249,"Segunda",1540,750,800,775,775),
    @(275,44265,"Primera",3400,800,900,850,850),
    @(276,44265,"Segunda",1740,700,750,725,725),
    @(277,44343,"Primera",2440,650,700,675,675),
    @(278,44343,"Segunda",1560,500,550,525,525),
    @(279,44201,"Primera",2000,650,700,675,675),
    @(280,44201,"Segunda",1300,550,600,575,575),
    @(281,44280,"Primera",2500,700,750,725,725),
    @(282,44280,"Segunda",1440,600,650,625,625),
    @(283,44447,"Primera",3300,650,700,675,675),
    @(284,44447,"Segunda",1600,550,600,575,575),
    @(285,44270,"Primera",2800,850,900,875,875),
    @(286,44270,"Segunda",1540,750,800,775,775),
    @(287,44260,"Primera",3000,850,900,875,875),
    @(288,44260,"Segunda",1680,750,800,775,775),
    @(289,44267,"Primera",2800,850,900,875,875),
    @(290,44267,"Segunda",1520,750,800,775,775),
    @(291,44312,"Primera",3000,650,700,675,675),
    @(292,44312,"Segunda",1600,550,600,575,575),
    @(293,44187,"Primera",2400,550,600,575,575),
    @(294,44187,"Segunda",1400,450,500,475,475),
    @(295,44390,"Primera",2200,600,700,650,650),
    @(296,44390,"Segunda",1400,500,550,525,525),
    @(297,44386,"Primera",3360,600,700,650,650),
    @(298,44386,"Segunda",1600,500,550,525,525),
    @(299,44308,"Primera",2500,650,700,675,675),
    @(300,44308,"Segunda",1480,550,600,575,575),
    @(301,44264,"Primera",2800,800,900,850,850),
    @(302,44264,"Segunda",1600,700,750,725,725),
    @(303,44463,"Primera",3600,600,700,650,650),
    @(304,44463,"Segunda",1800,500,550,525,525),
    @(305,44196,"Primera",2000,550,600,575,575),
    @(306,44196,"Segunda",1460,450,500,475,475),
    @(307,44301,"Primera",2500,650,700,675,675),
    @(308,44301,"Segunda",1400,550,600,575,575),
    @(309,44251,"Segunda",1700,750,800,775,775),
    @(310,44243,"Primera",2400,650,700,675,675),
    @(311,44243,"Segunda",1400,550,600,575,575),
    @(312,44252,"Segunda",1480,750,800,775,775),
    @(313,44166,"Primera",2000,450,500,475,475),
    @(314,44166,"Segunda",1400,350,400,375,375),
    @(315,44168,"Primera",2000,450,500,475,475),
    @(316,44168,"Segunda",1300,350,400,375,375),
    @(317,44369,"Primera",2300,600,700,650,650),
    @(318,44369,"Segunda",1400,500,550,525,525),
    @(319,44433,"Primera",3320,650,700,675,675),
    @(320,44433,"Segunda",1640,550,600,575,575),
    @(321,44221,"Primera",3000,650,700,675,675),
    @(322,44221,"Segunda",1600,550,600,575,575),
    @(323,44371,"Primera",2300,600,700,650,650),
    @(324,44371,"Segunda",1400,500,550,525,525),
    @(325,44316,"Primera",3200,650,700,675,675),
    @(326,44316,"Segunda",1680,500,550,525,525),
    @(327,44279,"Primera",3200,700,750,725,725),
    @(328,44279,"Segunda",1680,600,650,625,625),
    @(329,44397,"Primera",2400,700,800,750,750),
    @(330,44397,"Segunda",1400,500,600,550,550),
    @(331,44363,"Primera",3280,600,700,650,650),
    @(332,44363,"Segunda",1600,500,550,525,525),
    @(333,44277,"Primera",2800,700,750,725,725),
    @(334,44277,"Segunda",1540,600,650,625,625),
    @(335,44291,"Primera",2800,650,700,675,675),
    @(336,44291,"Segunda",1660,550,600,575,575),
    @(337,44273,"Primera",2500,750,800,775,775),
    @(338,44273,"Segunda",1440,650,700,675,675),
    @(339,44438,"Primera",3340,650,700,675,675),
    @(340,44438,"Segunda",1600,550,600,575,575),
    @(341,44372,"Primera",3320,600,700,650,650),
    @(342,44372,"Segunda",1600,500,550,525,525),
    @(343,44286,"Primera",3200,750,800,775,775),
    @(344,44286,"Segunda",1720,650,700,675,675),
    @(345,44209,"Primera",2700,650,700,675,675),
    @(346,44209,"Segunda",1600,500,550,525,525),
    @(347,44356,"Primera",3300,600,700,650,650),
    @(348,44356,"Segunda",1700,500,550,525,525),
    @(349,44160,"Primera",2700,450,500,475,475),
    @(350,44160,"Segunda",1600,350,400,375,375),
    @(351,44351,"Primera",3300,600,700,650,650),
    @(352,44351,"Segunda",1600,500,550,525,525),
    @(353,44365,"Primera",3300,600,700,650,650),
    @(354,44365,"Segunda",1600,500,550,525,525),
    @(355,44306,"Primera",2800,650,700,675,675),
    @(356,44306,"Segunda",1460,550,600,575,575),
    @(357,44215,"Primera",2400,600,650,625,625),
    @(358,44215,"Segunda",1360,500,550,525,525),
    @(359,44175,"Primera",2000,550,600,575,575),
    @(360,44175,"Segunda",1400,450,500,475,475),
    @(361,44357,"Primera",2400,600,700,650,650),
    @(362,44357,"Segunda",1400,500,550,525,525),
    @(363,44203,"Primera",2000,650,700,675,675),
    @(364,44203,"Segunda",1400,550,600,575,575),
    @(365,44162,"Primera",2800,450,500,475,475),
    @(366,44162,"Segunda",1600,350,400,375,375),
    @(367,44410,"Primera",3400,650,700,675,675),
    @(368,44410,"Segunda",1600,550,600,575,575),
    @(369,44411,"Primera",2300,650,700,675,675),
    @(370,44411,"Segunda",1400,550,600,575,575),
    @(371,44257,"Primera",2600,850,900,875,875),
    @(372,44257,"Segunda",1560,750,800,775,775),
    @(373,44244,"Primera",3000,750,800,775,775),
    @(374,44244,"Segunda",1600,650,700,675,675),
    @(375,44176,"Primera",2800,550,600,575,575),
    @(376,44176,"Segunda",1600,450,500,475,475),
    @(377,44239,"Primera",3000,650,700,675,675),
    @(378,44239,"Segunda",1600,550,600,575,575),
    @(379,44376,"Primera",2500,600,700,650,650),
    @(380,44376,"Segunda",1400,500,550,525,525),
    @(381,44292,"Primera",2400,650,700,675,675),
    @(382,44292,"Segunda",1480,550,600,575,575),
    @(383,44358,"Primera",3340,600,700,650,650),
    @(384,44358,"Segunda",1600,500,550,525,525),
    @(385,44211,"Primera",2700,650,700,675,675),
    @(386,44211,"Segunda",1600,550,600,575,575),
    @(387,44425,"Primera",2200,650,700,675,675),
    @(388,44425,"Segunda",1400,550,600,575,575),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D Fecha
    $ws.Cells.Item($r, 9).Value  = $row[2]   # I Calidad
    $ws.Cells.Item($r, 10).Value = $row[3]   # J Volumen
    $ws.Cells.Item($r, 11).Value = $row[4]   # K Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[5]   # L Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[6]   # M Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[7]   # P Precio $/Kg
}

# New rows 389-390 (full rows, carrying the two date blocks pushed off the end)
$newRow389 = @{
    1 = 8
    2 = "Terminal La Palmera de La Serena"
    3 = "Coquimbo"
    4 = 44323
    5 = 4
    6 = 100112008
    7 = "Coliflor"
    8 = "Sin especificar"
    9 = "Primera"
    10 = 3280
    11 = 650
    12 = 700
    13 = 675
    14 = "`$/unidad"
    15 = "Provincia del Elquí"
    16 = 675
    17 = 1
    18 = "Hortaliza"
}
foreach ($col in $newRow389.Keys) {
    $ws.Cells.Item(389, $col).Value = $newRow389[$col]
}
$ws.Cells.Item(389, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$newRow390 = @{
    1 = 8
    2 = "Terminal La Palmera de La Serena"
    3 = "Coquimbo"
    4 = 44323
    5 = 4
    6 = 100112008
    7 = "Coliflor"
    8 = "Sin especificar"
    9 = "Segunda"
    10 = 1660
    11 = 500
    12 = 550
    13 = 525
    14 = "`$/unidad"
    15 = "Provincia del Elquí"
    16 = 525
    17 = 1
    18 = "Hortaliza"
}
foreach ($col in $newRow390.Keys) {
    $ws.Cells.Item(390, $col).Value = $newRow390[$col]
}
$ws.Cells.Item(390, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
